# Fix the mis-spelled header text first (shared string used by the
# "MODEL_CONDITION" header cell, currently E1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "MODELCONDITION"

# Delete column A entirely - this shifts B->A, C->B, D->C, E->D, F->E,
# matching the diff which drops the old "A" (1 / 15) column and moves
# everything one column to the left.
$ws.Columns("A").Delete()
